$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Worksheet, $Address, $NewValue)
    $range = $Worksheet.Range($Address)
    $range.NumberFormat = "@"
    $range.Value = $NewValue
    $range.Style = "Normal"
}

Set-TextValue $ws "D2" "46.268.43"
Set-TextValue $ws "E2" "  -0.51%  "
Set-TextValue $ws "D3" "2.608.82"
Set-TextValue $ws "E3" "  +2.68%  "
Set-TextValue $ws "E4" "  +0.10%  "
Set-TextValue $ws "D5" "306.83"
Set-TextValue $ws "E5" "  -0.08%  "
Set-TextValue $ws "D6" "99.42"
Set-TextValue $ws "E6" "  -4.38%  "
Set-TextValue $ws "D7" "0.601"
Set-TextValue $ws "E7" "  -1.42%  "
Set-TextValue $ws "E8" "  -0.14%  "
Set-TextValue $ws "D9" "0.578"
Set-TextValue $ws "E9" "  +0.83%  "
Set-TextValue $ws "D10" "39.39"
Set-TextValue $ws "E10" "  +0.30%  "
Set-TextValue $ws "D11" "0.0842"
Set-TextValue $ws "E11" "  +1.14%  "
Set-TextValue $ws "D12" "54.09"
Set-TextValue $ws "E12" "  -0.91%  "
Set-TextValue $ws "D13" "8.10"
Set-TextValue $ws "E13" "  +1.47%  "
Set-TextValue $ws "D14" "3.007.63"
Set-TextValue $ws "E14" "  +2.82%  "
Set-TextValue $ws "E15" "  +0.64%  "
Set-TextValue $ws "D16" "2.609.11"
Set-TextValue $ws "E16" "  +1.75%  "
Set-TextValue $ws "D17" "0.917"
Set-TextValue $ws "E17" "  +2.49%  "
Set-TextValue $ws "D18" "14.92"
Set-TextValue $ws "E18" "  -0.47%  "
Set-TextValue $ws "D19" "46.446.03"
Set-TextValue $ws "E19" "  -0.29%  "
Set-TextValue $ws "E20" "  +0.60%  "
Set-TextValue $ws "D21" "12.89"
Set-TextValue $ws "E21" "  -7.56%  "
Set-TextValue $ws "D22" "6.70"
Set-TextValue $ws "E22" "  +0.89%  "
Set-TextValue $ws "D23" "71.31"
Set-TextValue $ws "E23" "  +1.58%  "
Set-TextValue $ws "D24" "271.96"
Set-TextValue $ws "E24" "  +6.59%  "
Set-TextValue $ws "D25" "3.03"
Set-TextValue $ws "E25" "  +1.40%  "
Set-TextValue $ws "E26" "  +1.43%  "
Set-TextValue $ws "D27" "29.00"
Set-TextValue $ws "E27" "  +20.37%  "
Set-TextValue $ws "D28" "1.00"
Set-TextValue $ws "E28" "  -0.03%  "
Set-TextValue $ws "E29" "  -0.69%  "
Set-TextValue $ws "D30" "10.56"
Set-TextValue $ws "E30" "  +0.94%  "
Set-TextValue $ws "D31" "38.54"
Set-TextValue $ws "E31" "  -8.69%  "
Set-TextValue $ws "D32" "2.20"
Set-TextValue $ws "E32" "  -2.99%  "
Set-TextValue $ws "D33" "6.32"
Set-TextValue $ws "E33" "  +4.78%  "
Set-TextValue $ws "D34" "3.64"
Set-TextValue $ws "E34" "  -5.33%  "
Set-TextValue $ws "E35" "  -1.95%  "
Set-TextValue $ws "D36" "2.23"
Set-TextValue $ws "E36" "  +1.48%  "
Set-TextValue $ws "D37" "0.0834"
Set-TextValue $ws "E37" "  -1.33%  "
Set-TextValue $ws "D38" "150.93"
Set-TextValue $ws "E38" "  +0.48%  "
Set-TextValue $ws "E39" "  +3.64%  "
Set-TextValue $ws "E40" "  +1.15%  "
Set-TextValue $ws "D41" "23.11"
Set-TextValue $ws "E41" "  +31.03%  "
Set-TextValue $ws "D42" "15.84"
Set-TextValue $ws "E42" "  -3.92%  "
Set-TextValue $ws "B43" "NEARProtocol"
Set-TextValue $ws "C43" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws "D43" "3.62"
Set-TextValue $ws "E43" "  +1.44%  "
Set-TextValue $ws "B44" "VeChain"
Set-TextValue $ws "C44" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws "D44" "0.0329"
Set-TextValue $ws "E44" "  +0.21%  "
Set-TextValue $ws "D45" "4.05"
Set-TextValue $ws "E45" "  -4.83%  "
Set-TextValue $ws "D46" "2.112.15"
Set-TextValue $ws "E46" "  +4.96%  "
Set-TextValue $ws "D47" "0.998"
Set-TextValue $ws "E47" "  -0.07%  "
Set-TextValue $ws "D48" "93.30"
Set-TextValue $ws "E48" "  -1.02%  "
Set-TextValue $ws "E49" "  +6.35%  "
Set-TextValue $ws "D50" "1.77"
Set-TextValue $ws "E50" "  -5.90%  "
Set-TextValue $ws "D51" "108.50"
Set-TextValue $ws "E51" "  +1.21%  "
